# Incorporate recent 3.0 updates
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "About" sheet updates
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# Source year and citation details (new EIA source document)
$wsAbout.Range("B4").Value = 2014
$wsAbout.Range("B5").Value = "Price Elasticities for Energy Use in Buildings of the United States"
$wsAbout.Range("B6").Value = "https://www.eia.gov/analysis/studies/buildings/energyuse/pdf/price_elasticities.pdf"
$wsAbout.Range("B7").Value = "Appendix"

# Note text now references 3-year (rather than 1-year) short-run elasticities
$wsAbout.Range("A10").Value = "We use same-price, long-run elasticities minus the 3-year short-run elasticities."
$wsAbout.Range("A11").Value = "We calculate it this way because we assume that 3-year elasticities primarily reflect behavior"
$wsAbout.Range("A14").Value = "all timescales.  So, the portion of the long-run elasticitiy represented by the 3-year elasticity"

# ---------------------------------------------------------------------
# "EIA Table 1" sheet updates (Residential block, rows 7-9)
# ---------------------------------------------------------------------
$wsEia = $wb.Worksheets.Item("EIA Table 1")

$wsEia.Range("B7").Value = -0.12
$wsEia.Range("C7").Value = -0.21
$wsEia.Range("D7").Value = -0.25
$wsEia.Range("E7").Value = -0.28000000000000003
$wsEia.Range("F7").Value = 0
$wsEia.Range("G7").Value = 0

$wsEia.Range("B8").Value = -0.07
$wsEia.Range("C8").Value = -0.13
$wsEia.Range("D8").Value = -0.15
$wsEia.Range("E8").Value = 0.03
$wsEia.Range("F8").Value = -0.21
$wsEia.Range("G8").Value = 0

$wsEia.Range("B9").Value = -0.07
$wsEia.Range("C9").Value = -0.12
$wsEia.Range("D9").Value = -0.14000000000000001
$wsEia.Range("E9").Value = 0
$wsEia.Range("F9").Value = 0
$wsEia.Range("G9").Value = -0.22

# "EIA Table 1" sheet updates (Commercial block, rows 14-16)
$wsEia.Range("B14").Value = -0.11
$wsEia.Range("C14").Value = -0.18
$wsEia.Range("D14").Value = -0.22
$wsEia.Range("E14").Value = -0.33
$wsEia.Range("F14").Value = 0.09
$wsEia.Range("G14").Value = 0

$wsEia.Range("B15").Value = -0.15
$wsEia.Range("C15").Value = -0.25
$wsEia.Range("D15").Value = -0.3
$wsEia.Range("E15").Value = 0.15
$wsEia.Range("F15").Value = -0.58
$wsEia.Range("G15").Value = 0.02

$wsEia.Range("B16").Value = -0.14000000000000001
$wsEia.Range("C16").Value = -0.24
$wsEia.Range("D16").Value = -0.28999999999999998
$wsEia.Range("E16").Value = 0
$wsEia.Range("F16").Value = 0.05
$wsEia.Range("G16").Value = -0.42

# ---------------------------------------------------------------------
# "EoCEDwEC" sheet updates - formulas now diff against column D
# (the 3-year elasticity column) rather than column B (the 1-year one)
# ---------------------------------------------------------------------
$wsEoc = $wb.Worksheets.Item("EoCEDwEC")

$wsEoc.Range("B2").Formula = "='EIA Table 1'!E7-'EIA Table 1'!D7"
$wsEoc.Range("D2").Formula = "='EIA Table 1'!E14-'EIA Table 1'!D14"
$wsEoc.Range("B4").Formula = "='EIA Table 1'!F8-'EIA Table 1'!D8"
$wsEoc.Range("D4").Formula = "='EIA Table 1'!F15-'EIA Table 1'!D15"
$wsEoc.Range("B5").Formula = "='EIA Table 1'!G9-'EIA Table 1'!D9"
$wsEoc.Range("D5").Formula = "='EIA Table 1'!G16-'EIA Table 1'!D16"

# ---------------------------------------------------------------------
# Sheet view / selection state
# ---------------------------------------------------------------------
$wsAbout.Range("A27").Select()
$wsEia.Range("E17").Select()
$wsEoc.Activate()
$wsEoc.Range("H29").Select()
